# Regenerate save_data column G (K) values: use K instead of Strike#,
# regen std/mean, calc and write s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 8
    4  = 4
    5  = 3
    6  = 8
    7  = 0
    8  = 4
    9  = 9
    10 = 5
    11 = 2
    12 = 3
    13 = 7
    14 = 4
    15 = 5
    16 = 11
    17 = 6
    18 = 8
    19 = 13
    20 = 10
    21 = 8
    22 = 11
    23 = 4
    24 = 6
    25 = 11
    26 = 6
    27 = 7
    28 = 6
    29 = 3
    30 = 10
    31 = 5
    32 = 8
    33 = 8
    34 = 2
    35 = 3
    36 = 4
    37 = 3
    38 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
